# Add summation and verify total net payment
# - Duplicate Sheet1's full data onto a new Sheet2 (kept as a reference/backup copy)
# - Add a "Curr" column (values "SGD") to both sheets
# - Trim Sheet1 down to the header + first two data rows (rows 1-3), since the
#   remaining rows now live on Sheet2

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Create Sheet2 right after Sheet1 and copy all of Sheet1's data into it.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"
$ws1.UsedRange.Copy($ws2.Range("A1"))

# New currency column header + values on Sheet2 (rows 2-12 are all "SGD").
$ws2.Range("I1").Value = "Curr"
for ($r = 2; $r -le 12; $r++) {
    $ws2.Cells.Item($r, 9).Value = "SGD"
}
$ws2.Cells.Select()

# New currency column header + values on Sheet1 (rows 2-3 only survive the trim below).
$ws1.Range("I1").Value = "Curr"
$ws1.Range("I2").Value = "SGD"
$ws1.Range("I3").Value = "SGD"

# Drop rows 4-12 from Sheet1 - that data now lives on Sheet2.
$ws1.Range("A4:I12").Select()
$ws1.Range("A4:I12").Delete()

$ws1.Activate()
$ws1.Range("A4:I12").Select()
